$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Logistic Regression) - slightly adjusted C2/E2 values
$ws.Range("C2").Value = 0.7446518136695224
$ws.Range("E2").Value = 0.737719644064865

# Row 3 becomes LightGBM with the values that used to belong to the (now removed) LightGBM row
$ws.Range("A3").Value = "LightGBM"
$ws.Range("B3").Value = 0.7730027548209366
$ws.Range("C3").Value = 0.7721349875663966
$ws.Range("D3").Value = 0.7730027548209366
$ws.Range("E3").Value = 0.7703176406920125

# Remove the now-unused rows (Support Vector Classifier, CART, Random Forest, LightGBM(old), XGBoost)
$ws.Range("A4:E8").Delete()
